$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.215.43"
$ws.Range("E2").Value = "  +0.89%  "
$ws.Range("D3").Value = "1.852.99"
$ws.Range("E3").Value = "  +1.52%  "
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("D5").Value = "'313.17"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("E6").Value = "  -0.44%  "
$ws.Range("D7").Value = "'0.4632"
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").Value = "'0.3717"
$ws.Range("E8").Value = "  +0.37%  "
$ws.Range("D9").Value = "'0.07283"
$ws.Range("E9").Value = "  -0.68%  "
$ws.Range("D10").Value = "'0.8866"
$ws.Range("E10").Value = "  +1.43%  "
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").Value = "'20.04"
$ws.Range("E11").Value = "  +1.25%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.07805"
$ws.Range("E12").Value = "  -1.54%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.806.29"
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("D14").Value = "'5.372"
$ws.Range("E14").Value = "  +0.65%  "
$ws.Range("D15").Value = "'6.512"
$ws.Range("E15").Value = "  -0.45%  "
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("D18").Value = "'0.000008913"
$ws.Range("E18").Value = "  +0.60%  "
$ws.Range("E19").Value = "  -0.33%  "
$ws.Range("D20").Value = "'14.69"
$ws.Range("D21").Value = "27.252.49"
$ws.Range("E21").Value = "  +0.91%  "
$ws.Range("D22").Value = "'5.059"
$ws.Range("E22").Value = "  -0.79%  "
$ws.Range("E23").Value = "  -0.27%  "
$ws.Range("D24").Value = "2.167.14"
$ws.Range("E24").Value = "  +4.79%  "
$ws.Range("D25").Value = "'1.950"
$ws.Range("E25").Value = "  +5.61%  "
$ws.Range("D26").Value = "'151.86"
$ws.Range("D27").Value = "'18.44"
$ws.Range("E27").Value = "  +0.20%  "
$ws.Range("D28").Value = "'2.042"
$ws.Range("E28").Value = "  +0.37%  "
$ws.Range("D29").Value = "'115.64"
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").Value = "'5.056"
$ws.Range("E30").Value = "  -1.46%  "
$ws.Range("D31").Value = "'0.08800"
$ws.Range("E31").Value = "  -1.11%  "
$ws.Range("D32").Value = "'3.138"
$ws.Range("E32").Value = "  +5.89%  "
$ws.Range("D33").Value = "'0.7651"
$ws.Range("E33").Value = "  +5.28%  "
$ws.Range("D34").Value = "'1.170"
$ws.Range("E34").Value = "  +3.55%  "
$ws.Range("E35").Value = "  +1.47%  "
$ws.Range("D36").Value = "'2.748"
$ws.Range("E36").Value = "  +11.34%  "
$ws.Range("D37").Value = "'1.087"
$ws.Range("E37").Value = "  +1.79%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "'0.01939"
$ws.Range("E38").Value = "  -0.40%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.05228"
$ws.Range("E39").Value = "  +0.08%  "
$ws.Range("D40").Value = "'2.930"
$ws.Range("E40").Value = "  -0.25%  "
$ws.Range("D41").Value = "'7.073"
$ws.Range("E41").Value = "  -0.36%  "
$ws.Range("D42").Value = "'0.5101"
$ws.Range("E42").Value = "  -1.02%  "
$ws.Range("D43").Value = "'0.1628"
$ws.Range("E43").Value = "  +0.57%  "
$ws.Range("D44").Value = "'8.397"
$ws.Range("E44").Value = "  +2.75%  "
$ws.Range("D45").Value = "'0.4783"
$ws.Range("E45").Value = "  -1.16%  "
$ws.Range("E46").Value = "  +1.23%  "
$ws.Range("E47").Value = "  -0.49%  "
$ws.Range("D48").Value = "'102.99"
$ws.Range("E48").Value = "  +0.17%  "
$ws.Range("D49").Value = "'1.638"
$ws.Range("E49").Value = "  +0.40%  "
$ws.Range("D50").Value = "'0.06212"
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("D51").Value = "'65.40"
$ws.Range("E51").Value = "  +1.20%  "
